# Apply cryptos list price/volume update (commit: "Updated cryptos list on Wed Jul 26 20:28:42 UTC 2023 with GitHub Actions")
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "29.464.99"
$ws.Range("E2").Value = "  +0.76%  "
$ws.Range("D3").Value = "1.877.26"
$ws.Range("E3").Value = "  +0.94%  "
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "0.9996"
$ws.Range("E4").Value = "  -0.17%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "0.7126"
$ws.Range("E5").Value = "  +1.67%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "238.94"
$ws.Range("E6").Value = "  +0.67%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.9989"
$ws.Range("E7").Value = "  -0.22%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.07891"
$ws.Range("E8").Value = "  -3.11%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.3081"
$ws.Range("E9").Value = "  +1.92%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "25.39"
$ws.Range("E10").Value = "  +9.50%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.08200"
$ws.Range("E11").Value = "  +0.53%  "
$ws.Range("B12").Value = "WrappedEther"
$ws.Range("C12").Value = "https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"
$ws.Range("D12").Value = "1.867.91"
$ws.Range("E12").Value = "  +0.31%  "
$ws.Range("B13").Value = "Polkadot"
$ws.Range("C13").Value = "https://coinranking.com/coin/25W7FG7om+polkadot-dot"
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "5.271"
$ws.Range("E13").Value = "  +2.33%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "0.7254"
$ws.Range("E14").Value = "  +3.06%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "89.45"
$ws.Range("E15").Value = "  +0.51%  "
$ws.Range("D16").Value = "29.466.67"
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "5.840"
$ws.Range("E17").Value = "  +1.31%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "0.000007856"
$ws.Range("E18").Value = "  +0.32%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "241.78"
$ws.Range("E19").Value = "  +2.76%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "13.37"
$ws.Range("E20").Value = "  +0.48%  "
$ws.Range("D21").Value = "2.120.91"
$ws.Range("E21").Value = "  +0.58%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "0.9991"
$ws.Range("E22").Value = "  -0.16%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "0.9998"
$ws.Range("E23").Value = "  -0.21%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "7.783"
$ws.Range("E24").Value = "  +5.13%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "163.35"
$ws.Range("E25").Value = "  +1.19%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "0.1471"
$ws.Range("E26").Value = "  +1.97%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "8.981"
$ws.Range("E27").Value = "  +0.48%  "
$ws.Range("E28").Value = "  +0.86%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "1.953"
$ws.Range("E29").Value = "  -0.18%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "1.358"
$ws.Range("E30").Value = "  -5.18%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "1.483"
$ws.Range("E31").Value = "  +0.29%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "4.344"
$ws.Range("E32").Value = "  -1.22%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "4.106"
$ws.Range("E33").Value = "  +1.33%  "
$ws.Range("E34").Value = "  +1.27%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "1.194"
$ws.Range("E35").Value = "  +2.46%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.7223"
$ws.Range("E36").Value = "  +2.33%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "1.001"
$ws.Range("E37").Value = "  +0.29%  "
$ws.Range("E38").Value = "  -0.20%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.01857"
$ws.Range("E39").Value = "  +1.24%  "
$ws.Range("E40").Value = "  -0.76%  "
$ws.Range("D41").Value = "1.176.76"
$ws.Range("E41").Value = "  +3.68%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.9126"
$ws.Range("E42").Value = "  -1.00%  "
$ws.Range("B43").Value = "Aave"
$ws.Range("C43").Value = "https://coinranking.com/coin/ixgUfzmLR+aave-aave"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "72.05"
$ws.Range("E43").Value = "  +2.82%  "
$ws.Range("B44").Value = "FraxShare"
$ws.Range("C44").Value = "https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "5.988"
$ws.Range("E44").Value = "  +1.77%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.4317"
$ws.Range("E45").Value = "  +1.48%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.9987"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "102.40"
$ws.Range("E47").Value = "  +0.24%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.5350"
$ws.Range("E48").Value = "  -1.53%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "1.776"
$ws.Range("E49").Value = "  +0.80%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "2.916"
$ws.Range("E50").Value = "  +6.54%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "9.230"
$ws.Range("E51").Value = "  +0.95%  "
